$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Insert the 4 new columns that hold the "Beelingua" fee/jatem data.
#    Target final layout: new columns land at P(16), Z(26), AM(39), AW(49).
#    Inserting left-to-right at these absolute positions works because each
#    index already accounts for the shift caused by the earlier insertions.
# ---------------------------------------------------------------------------
$ws.Columns.Item(16).Insert()   # new P: "(1) Beelingua Fee"
$ws.Columns.Item(26).Insert()   # new Z: "(1) Jatem Beelingua"
$ws.Columns.Item(39).Insert()   # new AM: "(2) Beelingua Fee"
$ws.Columns.Item(49).Insert()   # new AW: "(2) Jatem Beelingua"

# ---------------------------------------------------------------------------
# 2) Header row (row 1) - new header cells + copy neighbouring look & feel
# ---------------------------------------------------------------------------
$ws.Range("P1").Value = "(1)`nBeelingua Fee"
$ws.Range("P1").Interior.Color = $ws.Range("O1").Interior.Color
$ws.Range("P1").Font.Bold = $true
$ws.Range("P1").HorizontalAlignment = -4108
$ws.Range("P1").VerticalAlignment = -4108
$ws.Range("P1").WrapText = $true
$ws.Range("P1").NumberFormat = "_-* #,##0_-;\-* #,##0_-;_-* ""-""??_-;_-@_-"

$ws.Range("Z1").Value = "(1)`nJatem Beelingua"
$ws.Range("Z1").HorizontalAlignment = -4108
$ws.Range("Z1").VerticalAlignment = -4108
$ws.Range("Z1").WrapText = $true

$ws.Range("AM1").Value = "(2)`nBeelingua Fee"
$ws.Range("AM1").HorizontalAlignment = -4108
$ws.Range("AM1").VerticalAlignment = -4108
$ws.Range("AM1").NumberFormat = "_-* #,##0_-;\-* #,##0_-;_-* ""-""??_-;_-@_-"

$ws.Range("AW1").Value = "(2)`nJatem Beelingua"
$ws.Range("AW1").HorizontalAlignment = -4108
$ws.Range("AW1").VerticalAlignment = -4108
$ws.Range("AW1").WrapText = $true

# ---------------------------------------------------------------------------
# 3) Data row (row 2) - new value cells
# ---------------------------------------------------------------------------
$ws.Range("P2").Value = 1600000
$ws.Range("P2").NumberFormat = "_-* #,##0_-;\-* #,##0_-;_-* ""-""??_-;_-@_-"
$ws.Range("P2").HorizontalAlignment = -4108
$ws.Range("P2").VerticalAlignment = -4108

$ws.Range("Z2").Value = 44888
$ws.Range("Z2").NumberFormat = "dd/mm/yyyy"

$ws.Range("AM2").Value = 1600000
$ws.Range("AM2").NumberFormat = "_-* #,##0_-;\-* #,##0_-;_-* ""-""??_-;_-@_-"
$ws.Range("AM2").HorizontalAlignment = -4108
$ws.Range("AM2").VerticalAlignment = -4108

$ws.Range("AW2").Value = 44888
$ws.Range("AW2").NumberFormat = "dd/mm/yyyy"

# ---------------------------------------------------------------------------
# 4) Fix up the formulas that must now include the new columns
# ---------------------------------------------------------------------------
$ws.Range("Q2").Formula = "=SUM(H2:P2)"
$ws.Range("AN2").Formula = "=Y2+AE2+AG2+AK2+AL2+AM2"

# ---------------------------------------------------------------------------
# 5) Hyperlink anchor got shifted by the earlier Y->AA shift; make sure it
#    still points at the right cell.
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Item(1).Range = $ws.Range("AA2")

# ---------------------------------------------------------------------------
# 6) View tweaks described by the diff
# ---------------------------------------------------------------------------
$ws.Application.ActiveWindow.ScrollColumn = 12
$ws.Range("Z8").Select()
